$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

# Update the Status value from "active" to "draft"
$ws.Range("B6").Value = "draft"

# Update the Date value to the new timestamp
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# The two custom cell styles used in this table already declare
# vertical="top" + wrapText alignment, but the "apply alignment" flag
# was not turned on, so the wrap text never actually rendered.
# Turning WrapText on (it is already effectively top-aligned) makes
# Excel mark the alignment as applied for these cells.
$ws.Range("A1:B14").WrapText = $true
